$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.609.55'
$ws.Range("E2").Value = '  +0.63%  '
$ws.Range("D3").Value = '1.579.96'
$ws.Range("E3").Value = '  -0.58%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '213.45'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.19%  '
$ws.Range("E6").Value = '  -0.63%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '44.58'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.77%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '24.09'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.22%  '
$ws.Range("E10").Value = '  -1.79%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0591'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.42%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0892'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.70%  '
$ws.Range("D13").Value = '1.804.69'
$ws.Range("E13").Value = '  -0.67%  '
$ws.Range("D14").Value = '1.581.56'
$ws.Range("E14").Value = '  -0.74%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.71'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.06%  '
$ws.Range("D16").Value = '28.588.06'
$ws.Range("E16").Value = '  +0.48%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.520'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.09%  '
$ws.Range("E18").Value = '  -1.42%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '231.33'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.74%  '
$ws.Range("E20").Value = '  -0.94%  '
$ws.Range("E21").Value = '  -2.09%  '
$ws.Range("E22").Value = '  -0.01%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.91'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.61%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.18'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.64%  '
$ws.Range("E25").Value = '  +5.80%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '151.12'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.35%  '
$ws.Range("E27").Value = '  -1.15%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.45'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.69%  '
$ws.Range("E29").Value = '  -2.44%  '
$ws.Range("E30").Value = '  -0.04%  '
$ws.Range("E31").Value = '  +2.61%  '
$ws.Range("E32").Value = '  -2.01%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.20'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.36%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.11'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.93%  '
$ws.Range("D35").Value = '1.397.85'
$ws.Range("E35").Value = '  -0.18%  '
$ws.Range("E37").Value = '  -3.56%  '
$ws.Range("E38").Value = '  +0.53%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.66'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +4.17%  '
$ws.Range("E40").Value = '  -0.33%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.522'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -3.50%  '
$ws.Range("E42").Value = '  +0.00%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.792'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.20%  '
$ws.Range("E44").Value = '  +1.86%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0467'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.19%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '5.46'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.61%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.960'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.23%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '63.32'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.27%  '
$ws.Range("D49").Value = '1.716.88'
$ws.Range("E49").Value = '  -0.48%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '86.49'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.88%  '
$ws.Range("D51").Value = '0.0₆0102'
$ws.Range("E51").Value = '  -2.18%  '
